$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A138").Value = "2023-12-09 11:01:42"
$ws.Range("B138").Value = 0.0004

$ws.Range("A139").Value = "2023-12-09 11:01:51"
$ws.Range("B139").Value = 0.0002
